$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.912.33"
$ws.Range("E2").Value = "'  +2.19%  "
$ws.Range("D3").Value = "'1.814.02"
$ws.Range("E3").Value = "'  +2.93%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("D5").Value = "'311.95"
$ws.Range("E5").Value = "'  +2.58%  "
$ws.Range("E6").Value = "'  +0.41%  "
$ws.Range("D7").Value = "'0.4294"
$ws.Range("E7").Value = "'  +0.86%  "
$ws.Range("D8").Value = "'0.3699"
$ws.Range("E8").Value = "'  +2.25%  "
$ws.Range("D9").Value = "'0.07238"
$ws.Range("E9").Value = "'  +2.60%  "
$ws.Range("D10").Value = "'0.8679"
$ws.Range("E10").Value = "'  +4.54%  "
$ws.Range("D11").Value = "'2.091.48"
$ws.Range("E11").Value = "'  +22.73%  "
$ws.Range("D12").Value = "'21.33"
$ws.Range("E12").Value = "'  +5.95%  "
$ws.Range("D13").Value = "'6.636"
$ws.Range("E13").Value = "'  +3.86%  "
$ws.Range("D14").Value = "'5.400"
$ws.Range("E14").Value = "'  +3.20%  "
$ws.Range("D15").Value = "'0.06933"
$ws.Range("E15").Value = "'  +2.30%  "
$ws.Range("D16").Value = "'80.83"
$ws.Range("E16").Value = "'  +2.14%  "
$ws.Range("D17").Value = "'1.005"
$ws.Range("E17").Value = "'  +0.13%  "
$ws.Range("D18").Value = "'0.000008840"
$ws.Range("E18").Value = "'  +2.62%  "
$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "'  +0.43%  "
$ws.Range("E20").Value = "'  +1.79%  "
$ws.Range("D21").Value = "'26.935.41"
$ws.Range("E21").Value = "'  +4.40%  "
$ws.Range("D22").Value = "'5.201"
$ws.Range("E22").Value = "'  +4.12%  "
$ws.Range("D23").Value = "'10.97"
$ws.Range("E23").Value = "'  -0.87%  "
$ws.Range("D24").Value = "'2.319.88"
$ws.Range("E24").Value = "'  +20.17%  "
$ws.Range("D25").Value = "'154.64"
$ws.Range("E25").Value = "'  +1.73%  "
$ws.Range("D26").Value = "'1.886"
$ws.Range("E26").Value = "'  -0.96%  "
$ws.Range("D27").Value = "'18.37"
$ws.Range("E27").Value = "'  +1.42%  "
$ws.Range("D28").Value = "'5.250"
$ws.Range("E28").Value = "'  +4.85%  "
$ws.Range("D29").Value = "'1.929"
$ws.Range("E29").Value = "'  +15.38%  "
$ws.Range("D30").Value = "'114.57"
$ws.Range("E30").Value = "'  -0.10%  "
$ws.Range("D31").Value = "'0.08955"
$ws.Range("E31").Value = "'  +0.89%  "
$ws.Range("B32").Value = "'ARBITRUM"
$ws.Range("C32").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").Value = "'1.168"
$ws.Range("E32").Value = "'  +4.74%  "
$ws.Range("B33").Value = "'ImmutableX"
$ws.Range("C33").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.7447"
$ws.Range("E33").Value = "'  +3.42%  "
$ws.Range("D34").Value = "'4.432"
$ws.Range("E34").Value = "'  +3.10%  "
$ws.Range("D35").Value = "'2.803"
$ws.Range("E35").Value = "'  +3.52%  "
$ws.Range("D36").Value = "'1.005"
$ws.Range("E36").Value = "'  +0.59%  "
$ws.Range("D37").Value = "'1.115"
$ws.Range("E37").Value = "'  +4.45%  "
$ws.Range("D38").Value = "'0.05240"
$ws.Range("E38").Value = "'  +2.99%  "
$ws.Range("D39").Value = "'0.01922"
$ws.Range("E39").Value = "'  +2.21%  "
$ws.Range("D40").Value = "'0.5103"
$ws.Range("E40").Value = "'  +4.34%  "
$ws.Range("D41").Value = "'2.753"
$ws.Range("E41").Value = "'  +11.65%  "
$ws.Range("E42").Value = "'  +3.50%  "
$ws.Range("D43").Value = "'6.511"
$ws.Range("E43").Value = "'  +5.44%  "
$ws.Range("D44").Value = "'8.315"
$ws.Range("E44").Value = "'  +4.31%  "
$ws.Range("D45").Value = "'107.36"
$ws.Range("E45").Value = "'  +2.73%  "
$ws.Range("D46").Value = "'10.45"
$ws.Range("E46").Value = "'  +4.48%  "
$ws.Range("D47").Value = "'1.005"
$ws.Range("E47").Value = "'  +0.46%  "
$ws.Range("D48").Value = "'0.4586"
$ws.Range("E48").Value = "'  +3.00%  "
$ws.Range("D49").Value = "'1.651"
$ws.Range("E49").Value = "'  +5.49%  "
$ws.Range("D50").Value = "'0.06276"
$ws.Range("D51").Value = "'1.811"
$ws.Range("E51").Value = "'  +5.92%  "
